$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.786669850349426
$ws.Range("B1").Value = 4.155766487121582
$ws.Range("C1").Value = 2.084121227264404
$ws.Range("D1").Value = 0.8923735618591309
$ws.Range("E1").Value = 0.481493353843689
